$d = $word.ActiveDocument

# --- Part 1: format the title "Test Closure Report" ---
$titlePara = $d.Paragraphs.First
$titleRange = $titlePara.Range
$titleRange.Font.Size = 16
$titleRange.Font.SizeBi = 16
$titleRange.Font.Underline = 6

# --- Part 2: split "Reset App State button overlap in Firefox" ---
$d.Content.Find.Execute("Reset App State button overlap in Firefox", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Reset App State button overlap in Firefox", 2)
